$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.128.17"
$ws.Range("E2").Value = "  +1.22%  "
$ws.Range("D3").Value = "3.545.49"
$ws.Range("E3").Value = "  +1.54%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "617.40"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.96"
$ws.Range("E6").Value = "  +3.76%  "
$ws.Range("D7").Value = "3.541.49"
$ws.Range("E7").Value = "  +1.43%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.488"
$ws.Range("E9").Value = "  +1.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.146"
$ws.Range("E10").Value = "  +5.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.29"
$ws.Range("E11").Value = "  +5.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.437"
$ws.Range("E12").Value = "  +3.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000223"
$ws.Range("E13").Value = "  +2.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.03"
$ws.Range("E14").Value = "  +5.19%  "
$ws.Range("D15").Value = "4.160.09"
$ws.Range("E15").Value = "  +1.95%  "
$ws.Range("D16").Value = "3.557.83"
$ws.Range("E16").Value = "  +2.00%  "
$ws.Range("D17").Value = "68.811.78"
$ws.Range("E17").Value = "  +2.37%  "
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.73"
$ws.Range("E19").Value = "  +5.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.85"
$ws.Range("E20").Value = "  +5.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.95"
$ws.Range("E21").Value = "  +10.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "455.36"
$ws.Range("E22").Value = "  +2.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.639"
$ws.Range("E23").Value = "  +2.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.30"
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000130"
$ws.Range("E25").Value = "  +2.36%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "3.703.91"
$ws.Range("E26").Value = "  +2.06%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.46"
$ws.Range("E27").Value = "  +3.33%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.98"
$ws.Range("E29").Value = "  +8.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.56"
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.68"
$ws.Range("E31").Value = "  +6.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.170"
$ws.Range("E32").Value = "  +3.22%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.40"
$ws.Range("E34").Value = "  +4.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.07"
$ws.Range("E35").Value = "  +1.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.91"
$ws.Range("E36").Value = "  +3.47%  "
$ws.Range("D37").Value = "3.551.37"
$ws.Range("E37").Value = "  +1.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.23"
$ws.Range("E38").Value = "  +3.21%  "
$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.34"
$ws.Range("E40").Value = "  +6.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "177.92"
$ws.Range("E41").Value = "  +4.21%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0914"
$ws.Range("E43").Value = "  +4.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.55"
$ws.Range("E44").Value = "  +2.23%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.896"
$ws.Range("E45").Value = "  +1.50%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "30.33"
$ws.Range("E46").Value = "  +14.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.33"
$ws.Range("E47").Value = "  +5.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.54"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.65"
$ws.Range("E49").Value = "  +4.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.76"
$ws.Range("E50").Value = "  +3.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.261"
$ws.Range("E51").Value = "  +7.27%  "
